$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab) to reflect the new "through" date
$ws.Name = "Through 2022-05-08"

# Update the header cell text in I1 (shared string "2022 (through 05-07)" -> "2022 (through 05-08)")
$ws.Range("I1").Value = "2022 (through 05-08)"

# Update the May total (I6): 24 -> 28
$ws.Range("I6").Value = 28

# Update the grand Total row (I14): 576 -> 580
$ws.Range("I14").Value = 580
